$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(286, 1).Value = "2023-12-13 17:15:14"
$ws.Cells.Item(286, 2).Value = 0.0016

$ws.Cells.Item(287, 1).Value = "2023-12-13 17:16:48"
$ws.Cells.Item(287, 2).Value = 0.005600000000000001

$ws.Cells.Item(288, 1).Value = "2023-12-13 17:17:08"
$ws.Cells.Item(288, 2).Value = 0.0008

$ws.Cells.Item(289, 1).Value = "2023-12-13 17:17:14"
$ws.Cells.Item(289, 2).Value = 0.0004
